$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Cells.Item(28, 8).Value = 3872.5715  # H28: 3619.1333 -> 3872.5715
$ws.Cells.Item(28, 9).Value = 942.1667  # I28: 817.7143 -> 942.1667
$ws.Cells.Item(28, 11).Value = 942.1667  # K28: 817.7143 -> 942.1667
$ws.Cells.Item(28, 13).Value = -457.1667  # M28: -332.7143 -> -457.1667

# Row 58
$ws.Cells.Item(58, 8).Value = 2925  # H58: 2700 -> 2925
$ws.Cells.Item(58, 10).Value = 2925  # J58: 2700 -> 2925
$ws.Cells.Item(58, 12).Value = 8775  # L58: 8100 -> 8775
$ws.Cells.Item(58, 14).Value = -9075  # N58: -8400 -> -9075

# Row 70
$ws.Cells.Item(70, 8).Value = 3138.4285  # H70: 2900.2727 -> 3138.4285
$ws.Cells.Item(70, 9).Value = 1965.7142  # I70: 1465.5714 -> 1965.7142
$ws.Cells.Item(70, 10).Value = 3724.7856  # J70: 3569.8 -> 3724.7856
$ws.Cells.Item(70, 11).Value = 5897.142599999999  # K70: 4396.7142 -> 5897.142599999999
$ws.Cells.Item(70, 12).Value = 11174.3568  # L70: 10709.4 -> 11174.3568
$ws.Cells.Item(70, 13).Value = -5627.142599999999  # M70: -4126.7142 -> -5627.142599999999
$ws.Cells.Item(70, 14).Value = -11714.3568  # N70: -11249.4 -> -11714.3568

# Row 73
$ws.Cells.Item(73, 8).Value = 3138.4285  # H73: 2900.2727 -> 3138.4285
$ws.Cells.Item(73, 9).Value = 1965.7142  # I73: 1465.5714 -> 1965.7142
$ws.Cells.Item(73, 10).Value = 3724.7856  # J73: 3569.8 -> 3724.7856
$ws.Cells.Item(73, 11).Value = 5897.142599999999  # K73: 4396.7142 -> 5897.142599999999
$ws.Cells.Item(73, 12).Value = 11174.3568  # L73: 10709.4 -> 11174.3568
$ws.Cells.Item(73, 13).Value = -4961.142599999999  # M73: -3460.7142 -> -4961.142599999999
$ws.Cells.Item(73, 14).Value = -13046.3568  # N73: -12581.4 -> -13046.3568

# Row 99
$ws.Cells.Item(99, 8).Value = 3500.0625  # H99: 3722.8 -> 3500.0625
$ws.Cells.Item(99, 9).Value = 1434.25  # I99: 1859.3334 -> 1434.25
$ws.Cells.Item(99, 11).Value = 4302.75  # K99: 5578.0002 -> 4302.75
$ws.Cells.Item(99, 13).Value = -2804.75  # M99: -4080.0002 -> -2804.75

# Row 106
$ws.Cells.Item(106, 8).Value = 4206.7144  # H106: 3968.375 -> 4206.7144
$ws.Cells.Item(106, 9).Value = 4407.8335  # I106: 4106.7144 -> 4407.8335
$ws.Cells.Item(106, 11).Value = 4407.8335  # K106: 4106.7144 -> 4407.8335
$ws.Cells.Item(106, 13).Value = -3776.8335  # M106: -3475.7144 -> -3776.8335

# Row 107
$ws.Cells.Item(107, 8).Value = 573.6316  # H107: 1043.25 -> 573.6316
$ws.Cells.Item(107, 9).Value = 512.94116  # I107: 1038.1111 -> 512.94116
$ws.Cells.Item(107, 11).Value = 512.94116  # K107: 1038.1111 -> 512.94116
$ws.Cells.Item(107, 13).Value = 1407.05884  # M107: 881.8888999999999 -> 1407.05884

$ws = $wb.Worksheets.Item("ARM")
# Row 27
$ws.Cells.Item(27, 8).Value = 0  # H27: 10000 -> 0
$ws.Cells.Item(27, 10).Value = 0  # J27: 10000 -> 0
$ws.Cells.Item(27, 12).ClearContents()  # L27: delete (was 10000)
$ws.Cells.Item(27, 14).Value = 0  # N27: -10368 -> 0

# Row 46
$ws.Cells.Item(46, 8).Value = 0  # H46: 3330 -> 0
$ws.Cells.Item(46, 9).Value = 0  # I46: 3500 -> 0
$ws.Cells.Item(46, 10).Value = 0  # J46: 2990 -> 0
$ws.Cells.Item(46, 11).Value = 0  # K46: 3500 -> 0
$ws.Cells.Item(46, 12).Value = 0  # L46: 2990 -> 0
$ws.Cells.Item(46, 13).ClearContents()  # M46: delete (was -3181)
$ws.Cells.Item(46, 14).Value = 0  # N46: -3628 -> 0

# Row 110
$ws.Cells.Item(110, 8).Value = 882.17645  # H110: 932.3125 -> 882.17645
$ws.Cells.Item(110, 9).Value = 642.0714  # I110: 685.3077 -> 642.0714
$ws.Cells.Item(110, 11).Value = 642.0714  # K110: 685.3077 -> 642.0714
$ws.Cells.Item(110, 13).Value = 1402.9286  # M110: 1359.6923 -> 1402.9286

# Row 126
$ws.Cells.Item(126, 8).Value = 6222.222  # H126: 5925 -> 6222.222
$ws.Cells.Item(126, 9).Value = 6222.222  # I126: 5925 -> 6222.222
$ws.Cells.Item(126, 11).Value = 18666.666  # K126: 17775 -> 18666.666
$ws.Cells.Item(126, 13).Value = -16196.666  # M126: -15305 -> -16196.666

$ws = $wb.Worksheets.Item("BSM")
# Row 7
$ws.Cells.Item(7, 8).Value = 10004130  # H7: 10004540 -> 10004130
$ws.Cells.Item(7, 9).Value = 9171829  # I7: 10005736 -> 9171829
$ws.Cells.Item(7, 10).Value = 13333333  # J7: 10001250 -> 13333333
$ws.Cells.Item(7, 11).Value = 9171829  # K7: 10005736 -> 9171829
$ws.Cells.Item(7, 12).Value = 13333333  # L7: 10001250 -> 13333333
$ws.Cells.Item(7, 13).Value = -9171716  # M7: -10005623 -> -9171716
$ws.Cells.Item(7, 14).Value = -13333559  # N7: -10001476 -> -13333559

# Row 20
$ws.Cells.Item(20, 8).Value = 2572  # H20: 2955.875 -> 2572
$ws.Cells.Item(20, 9).Value = 1466.3334  # I20: 2069 -> 1466.3334
$ws.Cells.Item(20, 10).Value = 3124.8333  # J20: 3488 -> 3124.8333
$ws.Cells.Item(20, 11).Value = 1466.3334  # K20: 2069 -> 1466.3334
$ws.Cells.Item(20, 12).Value = 3124.8333  # L20: 3488 -> 3124.8333
$ws.Cells.Item(20, 13).Value = -1219.3334  # M20: -1822 -> -1219.3334
$ws.Cells.Item(20, 14).Value = -3618.8333  # N20: -3982 -> -3618.8333

# Row 86
$ws.Cells.Item(86, 8).Value = 5174.36  # H86: 5179.4 -> 5174.36
$ws.Cells.Item(86, 9).Value = 4170.5  # I86: 4270.3076 -> 4170.5
$ws.Cells.Item(86, 10).Value = 6452  # J86: 6164.25 -> 6452
$ws.Cells.Item(86, 11).Value = 4170.5  # K86: 4270.3076 -> 4170.5
$ws.Cells.Item(86, 12).Value = 6452  # L86: 6164.25 -> 6452
$ws.Cells.Item(86, 13).Value = -3047.5  # M86: -3147.3076 -> -3047.5
$ws.Cells.Item(86, 14).Value = -8698  # N86: -8410.25 -> -8698

# Row 89
$ws.Cells.Item(89, 8).Value = 5174.36  # H89: 5179.4 -> 5174.36
$ws.Cells.Item(89, 9).Value = 4170.5  # I89: 4270.3076 -> 4170.5
$ws.Cells.Item(89, 10).Value = 6452  # J89: 6164.25 -> 6452
$ws.Cells.Item(89, 11).Value = 20852.5  # K89: 21351.538 -> 20852.5
$ws.Cells.Item(89, 12).Value = 32260  # L89: 30821.25 -> 32260
$ws.Cells.Item(89, 13).Value = -15236.5  # M89: -15735.538 -> -15236.5
$ws.Cells.Item(89, 14).Value = -43492  # N89: -42053.25 -> -43492

$ws = $wb.Worksheets.Item("CRP")
# Row 47
$ws.Cells.Item(47, 8).Value = 7688  # H47: 6854.6665 -> 7688
$ws.Cells.Item(47, 9).Value = 7688  # I47: 7032 -> 7688
$ws.Cells.Item(47, 10).Value = 0  # J47: 6500 -> 0
$ws.Cells.Item(47, 11).Value = 7688  # K47: 7032 -> 7688
$ws.Cells.Item(47, 12).Value = 0  # L47: 6500 -> 0
$ws.Cells.Item(47, 13).ClearContents()  # M47: delete (was -6466)
$ws.Cells.Item(47, 14).Value = -7122  # N47: -7632 -> -7122

# Row 62
$ws.Cells.Item(62, 8).Value = 3722.5  # H62: 4191.25 -> 3722.5
$ws.Cells.Item(62, 9).Value = 3722.5  # I62: 4191.25 -> 3722.5
$ws.Cells.Item(62, 11).Value = 3722.5  # K62: 4191.25 -> 3722.5
$ws.Cells.Item(62, 13).Value = -3098.5  # M62: -3567.25 -> -3098.5

# Row 65
$ws.Cells.Item(65, 8).Value = 3722.5  # H65: 4191.25 -> 3722.5
$ws.Cells.Item(65, 9).Value = 3722.5  # I65: 4191.25 -> 3722.5
$ws.Cells.Item(65, 11).Value = 18612.5  # K65: 20956.25 -> 18612.5
$ws.Cells.Item(65, 13).Value = -15492.5  # M65: -17836.25 -> -15492.5

# Row 105
$ws.Cells.Item(105, 8).Value = 2754.4614  # H105: 2761.3076 -> 2754.4614
$ws.Cells.Item(105, 9).Value = 2761.8  # I105: 2649.6667 -> 2761.8
$ws.Cells.Item(105, 10).Value = 2749.875  # J105: 2857 -> 2749.875
$ws.Cells.Item(105, 11).Value = 2761.8  # K105: 2649.6667 -> 2761.8
$ws.Cells.Item(105, 12).Value = 2749.875  # L105: 2857 -> 2749.875
$ws.Cells.Item(105, 13).Value = -1014.8  # M105: -902.6667000000002 -> -1014.8
$ws.Cells.Item(105, 14).Value = -6243.875  # N105: -6351 -> -6243.875

# Row 132
$ws.Cells.Item(132, 8).Value = 994.5  # H132: 994.6 -> 994.5
$ws.Cells.Item(132, 9).Value = 994.5  # I132: 994.6 -> 994.5
$ws.Cells.Item(132, 11).Value = 2983.5  # K132: 2983.8 -> 2983.5
$ws.Cells.Item(132, 13).Value = -453.5  # M132: -453.8000000000002 -> -453.5

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Cells.Item(5, 8).Value = 1224.6  # H5: 1305 -> 1224.6
$ws.Cells.Item(5, 9).Value = 957.7778  # I5: 1065.125 -> 957.7778
$ws.Cells.Item(5, 11).Value = 2873.3334  # K5: 3195.375 -> 2873.3334
$ws.Cells.Item(5, 13).Value = -2761.3334  # M5: -3083.375 -> -2761.3334

# Row 11
$ws.Cells.Item(11, 8).Value = 100  # H11: 185.71428 -> 100
$ws.Cells.Item(11, 10).Value = 0  # J11: 200 -> 0
$ws.Cells.Item(11, 12).Value = 0  # L11: 600 -> 0
$ws.Cells.Item(11, 14).ClearContents()  # N11: delete (was -880)

# Row 26
$ws.Cells.Item(26, 8).Value = 39.714287  # H26: 44 -> 39.714287
$ws.Cells.Item(26, 9).Value = 44.5  # I26: 60 -> 44.5
$ws.Cells.Item(26, 11).Value = 133.5  # K26: 180 -> 133.5
$ws.Cells.Item(26, 13).Value = 154.5  # M26: 108 -> 154.5

# Row 52
$ws.Cells.Item(52, 8).Value = 2600  # H52: 2650 -> 2600
$ws.Cells.Item(52, 10).Value = 2600  # J52: 2650 -> 2600
$ws.Cells.Item(52, 12).Value = 7800  # L52: 7950 -> 7800
$ws.Cells.Item(52, 14).Value = -8332  # N52: -8482 -> -8332

# Row 59
$ws.Cells.Item(59, 8).Value = 0  # H59: 837.5 -> 0
$ws.Cells.Item(59, 9).Value = 0  # I59: 837.5 -> 0
$ws.Cells.Item(59, 11).Value = 0  # K59: 2512.5 -> 0
$ws.Cells.Item(59, 13).ClearContents()  # M59: delete (was -1972.5)

# Row 80
$ws.Cells.Item(80, 8).Value = 4251.4814  # H80: 4261.1924 -> 4251.4814
$ws.Cells.Item(80, 9).Value = 4055.5  # I80: 4058.8235 -> 4055.5
$ws.Cells.Item(80, 11).Value = 12166.5  # K80: 12176.4705 -> 12166.5
$ws.Cells.Item(80, 13).Value = -11230.5  # M80: -11240.4705 -> -11230.5

# Row 83
$ws.Cells.Item(83, 8).Value = 4251.4814  # H83: 4261.1924 -> 4251.4814
$ws.Cells.Item(83, 9).Value = 4055.5  # I83: 4058.8235 -> 4055.5
$ws.Cells.Item(83, 11).Value = 36499.5  # K83: 36529.4115 -> 36499.5
$ws.Cells.Item(83, 13).Value = -31819.5  # M83: -31849.4115 -> -31819.5

# Row 114
$ws.Cells.Item(114, 8).Value = 1772.6875  # H114: 1539.2667 -> 1772.6875
$ws.Cells.Item(114, 9).Value = 1062.375  # I114: 506.5 -> 1062.375
$ws.Cells.Item(114, 10).Value = 2483  # J114: 2227.7778 -> 2483
$ws.Cells.Item(114, 11).Value = 3187.125  # K114: 1519.5 -> 3187.125
$ws.Cells.Item(114, 12).Value = 7449  # L114: 6683.3334 -> 7449
$ws.Cells.Item(114, 13).Value = 66.875  # M114: 1734.5 -> 66.875
$ws.Cells.Item(114, 14).Value = -13957  # N114: -13191.3334 -> -13957

# Row 134
$ws.Cells.Item(134, 8).Value = 1143.3334  # H134: 1211.3334 -> 1143.3334
$ws.Cells.Item(134, 9).Value = 1143.3334  # I134: 1211.3334 -> 1143.3334
$ws.Cells.Item(134, 11).Value = 3430.0002  # K134: 3634.0002 -> 3430.0002
$ws.Cells.Item(134, 13).Value = 1639.9998  # M134: 1435.9998 -> 1639.9998

# Row 135
$ws.Cells.Item(135, 8).Value = 1224.6  # H135: 1305 -> 1224.6
$ws.Cells.Item(135, 9).Value = 957.7778  # I135: 1065.125 -> 957.7778
$ws.Cells.Item(135, 11).Value = 8620.0002  # K135: 9586.125 -> 8620.0002
$ws.Cells.Item(135, 13).Value = -6085.0002  # M135: -7051.125 -> -6085.0002

$ws = $wb.Worksheets.Item("GSM")
# Row 31
$ws.Cells.Item(31, 8).Value = 1378.6  # H31: 466.66666 -> 1378.6
$ws.Cells.Item(31, 9).Value = 1378.6  # I31: 466.66666 -> 1378.6
$ws.Cells.Item(31, 11).Value = 1378.6  # K31: 466.66666 -> 1378.6
$ws.Cells.Item(31, 13).Value = -1086.6  # M31: -174.66666 -> -1086.6

# Row 37
$ws.Cells.Item(37, 8).Value = 1378.6  # H37: 466.66666 -> 1378.6
$ws.Cells.Item(37, 9).Value = 1378.6  # I37: 466.66666 -> 1378.6
$ws.Cells.Item(37, 11).Value = 1378.6  # K37: 466.66666 -> 1378.6
$ws.Cells.Item(37, 13).Value = -1101.6  # M37: -189.66666 -> -1101.6

# Row 70
$ws.Cells.Item(70, 8).Value = 1453.4  # H70: 1470.4 -> 1453.4
$ws.Cells.Item(70, 9).Value = 1453.4  # I70: 1470.4 -> 1453.4
$ws.Cells.Item(70, 11).Value = 1453.4  # K70: 1470.4 -> 1453.4
$ws.Cells.Item(70, 13).Value = -1183.4  # M70: -1200.4 -> -1183.4

# Row 73
$ws.Cells.Item(73, 8).Value = 1453.4  # H73: 1470.4 -> 1453.4
$ws.Cells.Item(73, 9).Value = 1453.4  # I73: 1470.4 -> 1453.4
$ws.Cells.Item(73, 11).Value = 1453.4  # K73: 1470.4 -> 1453.4
$ws.Cells.Item(73, 13).Value = -517.4000000000001  # M73: -534.4000000000001 -> -517.4000000000001

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Cells.Item(16, 8).Value = 4467.3335  # H16: 7777 -> 4467.3335
$ws.Cells.Item(16, 9).Value = 2812.5  # I16: 0 -> 2812.5
$ws.Cells.Item(16, 11).Value = 2812.5  # K16: 0 -> 2812.5
$ws.Cells.Item(16, 13).Value = -2642.5  # M16: None -> -2642.5

# Row 20
$ws.Cells.Item(20, 8).Value = 28781.611  # H20: 28762.334 -> 28781.611
$ws.Cells.Item(20, 9).Value = 1076  # I20: 1027.8125 -> 1076
$ws.Cells.Item(20, 10).Value = 48571.332  # J20: 50949.95 -> 48571.332
$ws.Cells.Item(20, 11).Value = 1076  # K20: 1027.8125 -> 1076
$ws.Cells.Item(20, 12).Value = 48571.332  # L20: 50949.95 -> 48571.332
$ws.Cells.Item(20, 13).Value = -850  # M20: -801.8125 -> -850
$ws.Cells.Item(20, 14).Value = -49023.332  # N20: -51401.95 -> -49023.332

# Row 55
$ws.Cells.Item(55, 8).Value = 1565.3334  # H55: 765.7143 -> 1565.3334
$ws.Cells.Item(55, 9).Value = 1565.3334  # I55: 814.5833 -> 1565.3334
$ws.Cells.Item(55, 10).Value = 0  # J55: 472.5 -> 0
$ws.Cells.Item(55, 11).Value = 1565.3334  # K55: 814.5833 -> 1565.3334
$ws.Cells.Item(55, 12).Value = 0  # L55: 472.5 -> 0
$ws.Cells.Item(55, 13).ClearContents()  # M55: delete (was -641.5833)
$ws.Cells.Item(55, 14).Value = -1392.3334  # N55: -818.5 -> -1392.3334

# Row 61
$ws.Cells.Item(61, 8).Value = 2954.7  # H61: 3020.5 -> 2954.7
$ws.Cells.Item(61, 9).Value = 841.1667  # I61: 886.4286 -> 841.1667
$ws.Cells.Item(61, 10).Value = 6125  # J61: 8000 -> 6125
$ws.Cells.Item(61, 11).Value = 841.1667  # K61: 886.4286 -> 841.1667
$ws.Cells.Item(61, 12).Value = 6125  # L61: 8000 -> 6125
$ws.Cells.Item(61, 13).Value = -639.1667  # M61: -684.4286 -> -639.1667
$ws.Cells.Item(61, 14).Value = -6529  # N61: -8404 -> -6529

# Row 113
$ws.Cells.Item(113, 8).Value = 2954.7  # H113: 3020.5 -> 2954.7
$ws.Cells.Item(113, 9).Value = 841.1667  # I113: 886.4286 -> 841.1667
$ws.Cells.Item(113, 10).Value = 6125  # J113: 8000 -> 6125
$ws.Cells.Item(113, 11).Value = 841.1667  # K113: 886.4286 -> 841.1667
$ws.Cells.Item(113, 12).Value = 6125  # L113: 8000 -> 6125
$ws.Cells.Item(113, 13).Value = 1328.8333  # M113: 1283.5714 -> 1328.8333
$ws.Cells.Item(113, 14).Value = -10465  # N113: -12340 -> -10465

$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Cells.Item(107, 8).Value = 950  # H107: 394.14285 -> 950
$ws.Cells.Item(107, 9).Value = 0  # I107: 301.5 -> 0
$ws.Cells.Item(107, 11).Value = 0  # K107: 904.5 -> 0
$ws.Cells.Item(107, 13).ClearContents()  # M107: delete (was 1015.5)

# Row 113
$ws.Cells.Item(113, 8).Value = 638.375  # H113: 595 -> 638.375
$ws.Cells.Item(113, 9).Value = 641.6  # I113: 576 -> 641.6
$ws.Cells.Item(113, 11).Value = 1924.8  # K113: 1728 -> 1924.8
$ws.Cells.Item(113, 13).Value = 245.1999999999998  # M113: 442 -> 245.1999999999998

# Row 132
$ws.Cells.Item(132, 8).Value = 1842.5  # H132: 1914.7222 -> 1842.5
$ws.Cells.Item(132, 9).Value = 1781.579  # I132: 1850.8823 -> 1781.579
$ws.Cells.Item(132, 11).Value = 5344.737  # K132: 5552.6469 -> 5344.737
$ws.Cells.Item(132, 13).Value = -2814.737  # M132: -3022.6469 -> -2814.737
